$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / N values) - updated for "meanEMG legmaxROM" subject columns
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - updated values
$ws.Range("B2").Value = 7.159403109186818
$ws.Range("C2").Value = 5.0175368920160865
$ws.Range("D2").Value = 5.9453661633681243
$ws.Range("E2").Value = 7.1594031091868082

# Row 3 (STR) - updated values; C3 removed, D3 newly populated
$ws.Range("B3").Value = 6.6663762365620309
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 5.4767460211276076
$ws.Range("E3").Value = 6.1832668890764779

# Selection now reflects the updated/edited data block
$ws.Range("B1:E3").Select()
